$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers "Age (in years)" / "Education" (Kyrgyz/Russian/English columns) are
# reworded to "By age (in years)" / "By education". Cells are touched column
# by column (A, then B, then C) so new shared-string entries land in the same
# order as the canonical file.
$ws.Range("A19").Value = "Жаш курагы боюнча (жылдарда)"
$ws.Range("A29").Value = "Билими боюнча"

$ws.Range("B19").Value = "По возрасту (в годах)"
$ws.Range("B29").Value = "По образованию"

$ws.Range("C19").Value = "By age (in years) "
$ws.Range("C29").Value = "By education"
